# Remise en forme du code
# - Nettoyage du code
# - Creation du code en 'XXXX'
# - utilisation de char au lieu d'int

$wb = $excel.ActiveWorkbook

# --- Typography sheet: switch the "Wildcard Characters" for Font_Voleurs (row 6)
#     from plain digits to digits + X + hex letters (char-based code instead of int-based)
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("G6").Value = "0123456789XABCDEF"

# --- Translation sheet: add the new "code" related text entries
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B11").Value = "SingleUseId7"
$wsTrans.Range("C11").Value = "Default"
$wsTrans.Range("D11").Value = "Left"
$wsTrans.Range("E11").Value = "LTR"
$wsTrans.Range("F11").Value = "Code 1 : A14C`nCode 2 : 5C68`nCode 3 : 4489"

$wsTrans.Range("B12").Value = "SingleUseId8"
$wsTrans.Range("C12").Value = "Font_Voleurs"
$wsTrans.Range("D12").Value = "Left"
$wsTrans.Range("E12").Value = "LTR"
$wsTrans.Range("F12").Value = "<value>"

$wsTrans.Range("B13").Value = "SingleUseId9"
$wsTrans.Range("C13").Value = "Font_Voleurs"
$wsTrans.Range("D13").Value = "Left"
$wsTrans.Range("E13").Value = "LTR"
$wsTrans.Range("F13").Value = "X"

$wsTrans.Range("B14").Value = "SingleUseId10"
$wsTrans.Range("C14").Value = "Font_Voleurs"
$wsTrans.Range("D14").Value = "Left"
$wsTrans.Range("E14").Value = "LTR"
$wsTrans.Range("F14").Value = "<value>"

$wsTrans.Range("B15").Value = "SingleUseId11"
$wsTrans.Range("C15").Value = "Font_Voleurs"
$wsTrans.Range("D15").Value = "Left"
$wsTrans.Range("E15").Value = "LTR"
$wsTrans.Range("F15").Value = "X"

$wsTrans.Range("B16").Value = "SingleUseId12"
$wsTrans.Range("C16").Value = "Font_Voleurs"
$wsTrans.Range("D16").Value = "Left"
$wsTrans.Range("E16").Value = "LTR"
$wsTrans.Range("F16").Value = "<value>"

$wsTrans.Range("B17").Value = "SingleUseId13"
$wsTrans.Range("C17").Value = "Font_Voleurs"
$wsTrans.Range("D17").Value = "Left"
$wsTrans.Range("E17").Value = "LTR"
$wsTrans.Range("F17").Value = "X"

$wsTrans.Range("B18").Value = "SingleUseId14"
$wsTrans.Range("C18").Value = "Font_Voleurs"
$wsTrans.Range("D18").Value = "Left"
$wsTrans.Range("E18").Value = "LTR"
$wsTrans.Range("F18").Value = "<value>"

$wsTrans.Range("B19").Value = "SingleUseId15"
$wsTrans.Range("C19").Value = "Font_Voleurs"
$wsTrans.Range("D19").Value = "Left"
$wsTrans.Range("E19").Value = "LTR"
$wsTrans.Range("F19").Value = "X"

# Keep the new rows at the default row height (avoid Excel's automatic
# row-height growth from the multi-line "Code 1 : ..." text in F11).
$wsTrans.Range("B11:F19").Rows.AutoFit()
